$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$srcStyle = $ws.Range("C2").Style
$ws.Range("B14").Style = $srcStyle
